$d = $word.ActiveDocument

# --- 1. Remove the two trailing spacer runs (bold, sz 54) at the end of
#        the first paragraph (right after the big rectangle drawing). ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
# The paragraph holds the drawing run (empty text) followed by two runs of
# plain spaces ("       " + "        " = 15 spaces) and then the paragraph
# mark. Delete everything except the paragraph mark itself.
if ($r1.End - $r1.Start -gt 1) {
    $spacer = $d.Range($r1.Start, $r1.End - 1)
    if ($spacer.Text -match '^\s*$') {
        $spacer.Delete()
    }
}

# --- 2. Bump the font size used for the following empty paragraph's mark
#        (pPr/rPr) from 2pt (sz 4) up to 12pt (sz 24 / szCs 24). ---
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.InsertBefore("x")
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.Font.Size = 12
$r2.Font.SizeBi = 12
$killRange = $d.Range($r2.Start, $r2.Start + 1)
$killRange.Delete()
